$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Row 6: change "Rejected" -> "Approved" in I6, and clear J6 ("Not appropriate")
$ws.Range("I6").Value = "Approved"
$ws.Range("J6").Value = $null

# Row 8: change "Rejected" -> "Approved" in I8, and clear J8 ("Not required")
$ws.Range("I8").Value = "Approved"
$ws.Range("J8").Value = $null

# Update the current selection to J6 (matches the diff's sheetView selection)
$ws.Range("J6").Select()
